$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, pushing the existing rows 57-126 down to 58-127
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new weekly price record
$ws.Cells.Item(57, 1).Value = 5
$ws.Cells.Item(57, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(57, 3).Value = "Maule"
$ws.Cells.Item(57, 4).Value = 44629
$ws.Cells.Item(57, 5).Value = 7
$ws.Cells.Item(57, 6).Value = 100112030
$ws.Cells.Item(57, 7).Value = "Poroto granado"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 300
$ws.Cells.Item(57, 11).Value = 23000
$ws.Cells.Item(57, 12).Value = 23000
$ws.Cells.Item(57, 13).Value = 23000
$ws.Cells.Item(57, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(57, 15).Value = "Región del Maule"
$ws.Cells.Item(57, 16).Value = 920
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
